$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow columns A:C from 39 to 36.5703125 (matches diff's width change)
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 36.5703125

# Add a new "2023" column (T) mirroring the existing "2022" column (S)
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("T4").Value = 2023

$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("T5").Value = 40

$excel.CutCopyMode = 0

# Reset the view: clear the stale selection / frozen left column from the old state
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
